$wb = $excel.ActiveWorkbook

# Clear the computed exercise values (formulas + cached results) in EJERCICIO 2,
# keeping the existing cell formatting/styles intact.
$ws2 = $wb.Worksheets.Item("Hoja2")
$ws2.Range("B2:E14").ClearContents()

# Rename the sheets. Renaming automatically updates the "Radianes" defined name
# reference from Hoja2!... to 'EJERCICIO 2'!...
$wb.Worksheets.Item("Hoja1").Name = "EJERCICIO 1"
$wb.Worksheets.Item("Hoja2").Name = "EJERCICIO 2"

# Restore the active-sheet/selection state recorded in the file: EJERCICIO 2's
# selection moves to C37:C38, and EJERCICIO 1 (the tab shown on open) ends up
# with F28 selected.
$ws2 = $wb.Worksheets.Item("EJERCICIO 2")
$ws2.Activate()
$ws2.Range("C37:C38").Select() | Out-Null

$ws1 = $wb.Worksheets.Item("EJERCICIO 1")
$ws1.Activate()
$ws1.Range("F28").Select() | Out-Null
